$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells F1:H1, matching the style of the existing headers (A1:E1)
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

$headerSrc = $ws.Range("E1")
$headerDst = $ws.Range("F1:H1")
$headerDst.Font.Bold = $headerSrc.Font.Bold
$headerDst.HorizontalAlignment = $headerSrc.HorizontalAlignment
$headerDst.VerticalAlignment = $headerSrc.VerticalAlignment
$headerDst.Borders.LineStyle = $headerSrc.Borders.LineStyle

# Fill boolean columns F, G, H for rows 2-21 with FALSE by default
for ($row = 2; $row -le 21; $row++) {
    $ws.Cells.Item($row, 6).Value = $false
    $ws.Cells.Item($row, 7).Value = $false
    $ws.Cells.Item($row, 8).Value = $false
}

# Row 10 has KNN_Outliers_MAD = TRUE
$ws.Cells.Item(10, 6).Value = $true

# Update the used range dimension to match new data extent
$ws.UsedRange | Out-Null
